$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.435.12'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.185.21'
$ws.Range("E3").Value = '  -2.39%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.24'
$ws.Range("E5").Value = '  +2.63%  '
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.17'
$ws.Range("E7").Value = '  -1.19%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  -5.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.06'
$ws.Range("E10").Value = '  -3.58%  '
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.75'
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = '2.513.75'
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.16'
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("D16").Value = '2.192.72'
$ws.Range("E16").Value = '  -1.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.768'
$ws.Range("E17").Value = '  -5.82%  '
$ws.Range("D18").Value = '42.376.62'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.86'
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.86'
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.49'
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("E23").Value = '  -12.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.10'
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.43'
$ws.Range("E26").Value = '  -5.10%  '
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("E28").Value = '  -4.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.13'
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.95'
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.00'
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0820'
$ws.Range("E33").Value = '  +2.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.14'
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("E36").Value = '  -3.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.21'
$ws.Range("E37").Value = '  -3.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0332'
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.03'
$ws.Range("E39").Value = '  -8.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.07'
$ws.Range("E40").Value = '  -3.67%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.56'
$ws.Range("E41").Value = '  +10.66%  '
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.16'
$ws.Range("E42").Value = '  -7.81%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.193'
$ws.Range("E43").Value = '  -2.96%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '58.55'
$ws.Range("E44").Value = '  -3.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.60'
$ws.Range("E45").Value = '  -4.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0969'
$ws.Range("E46").Value = '  -2.75%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.17'
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("B48").Value = 'WOONetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.458'
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("E49").Value = '  -2.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.13'
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("E51").Value = '  -0.81%  '
